$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "survey" sheet: insert a new "end screen" row right after the geopoint
# question (old row 5) and remove the old "end screen" row that used to sit
# at the very end of the sheet (old row 8). Net effect: the form now ends
# right after the GPS-coordinates prompt, before the two select_one
# questions (water_body_type / stand_flow), which get pushed down by one row.
# ---------------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

# Push rows 6-8 down to make room for the new "end screen" row.
$survey.Rows("6:6").Insert()

# Populate the newly-inserted row with the "end screen" marker.
$survey.Range("A6").Value = "end screen"

# The insert operation carries formatting into C6/E6/F6; clear them since
# the new row only has a value in column A.
$survey.Range("C6:F6").Clear()

# Remove the now-redundant "end screen" row that used to terminate the form
# (originally row 8, now pushed down to row 9).
$survey.Rows("9:9").Delete()

# Update the remembered selection to match the edited workbook.
$survey.Range("A7").Select()

# ---------------------------------------------------------------------------
# "settings" sheet: bump the form_version setting to reflect the new build.
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("settings")
$settings.Range("B3").Value = 20210304001

# Re-activate the sheet/cell that was active before our edits (selecting a
# range on another sheet switches the active tab as a side effect, but the
# workbook's active tab did not change in the target edit).
$settings.Activate()
$settings.Range("B4").Select()
